$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B and C (text values: coin names and links) - set directly
$ws.Range("B15").Value = "CoinExToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("B20").Value = "MCDex"
$ws.Range("C20").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("B22").Value = "ZBToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"

# Column D and E (numeric-looking / percent-looking strings) - force text format per-cell to preserve exact string content
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "293.63"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-4.84%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "40.23"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-2.79%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.043"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-3.21%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07384"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-3.96%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.310"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-0.10%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.551"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-5.60%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9233"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.82%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1185"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-4.08%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1760"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-3.41%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08692"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-4.84%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04175"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-1.10%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.1055"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.44%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001254"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.35%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.03810"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-5.32%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005810"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.92%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.378"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.85%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.399"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-1.28%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3295"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-1.21%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.586"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "2.22%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1345"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-3.95%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2810"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.34%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001284"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.42%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.003644"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-11.11%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001294"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.69%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0003736"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-95.02%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02307"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "-9.61%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05020"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-5.92%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007696"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-2.03%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.004416"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "136.95%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-2.90%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007410"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "10.85%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.006990"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-5.41%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3188"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "4.13%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006472"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-3.88%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000753"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.17%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "29.56%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.004216"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "35.83%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002108"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.17%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002007"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.17%"

Write-Host "Applied cryptos.xlsx update"
